# Fruta / hortaliza, semanal
# Insert the latest weekly price record as a new row right after the header
# block of existing rows (new row 16), pushing the previously-existing rows
# 16..30 down to 17..31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16 (shifts rows 16:30 down to 17:31).
$ws.Rows(16).Insert()

# Populate the new row 16 with this week's record.
$ws.Range("A16").Value = 8
$ws.Range("B16").Value = "Terminal La Palmera de La Serena"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44665
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100101
$ws.Range("H16").Value = "Berries"
$ws.Range("I16").Value = 100101001
$ws.Range("J16").Value = "Arándano (blue)"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 160
$ws.Range("N16").Value = 6500
$ws.Range("O16").Value = 7000
$ws.Range("P16").Value = 6750
$ws.Range("Q16").Value = "`$/bandeja 2 kilos"
$ws.Range("R16").Value = "Provincia de Linares"
$ws.Range("S16").Value = 3375
$ws.Range("T16").Value = 2
